$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (the "Duable Face" entry); remaining rows shift up by one.
$ws.Rows.Item(5).Delete()
